$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is purely numeric-looking text must be written with a leading
# apostrophe so Excel keeps them as text (matching the original inlineStr cells)
# instead of silently converting them to numbers; the style is then reset back to
# "Normal" so no stray number-format / quote-prefix style sticks to the cell.

$ws.Range('D2').Value = '70.752.26'
$ws.Range('E2').Value = '  +7.59%  '
$ws.Range('D3').Value = '3.637.48'
$ws.Range('E3').Value = '  +7.48%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'594.95"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.83%  '
$ws.Range('D6').Value = "'192.09"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.13%  '
$ws.Range('D7').Value = "'0.649"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.00%  '
$ws.Range('D8').Value = '3.617.07'
$ws.Range('E8').Value = '  +7.07%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = "'0.181"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.45%  '
$ws.Range('D11').Value = "'0.664"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.89%  '
$ws.Range('D12').Value = "'58.00"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.96%  '
$ws.Range('D13').Value = "'0.0000296"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.70%  '
$ws.Range('E14').Value = '  +5.90%  '
$ws.Range('D15').Value = '4.213.99'
$ws.Range('E15').Value = '  +7.36%  '
$ws.Range('D16').Value = '3.633.07'
$ws.Range('E16').Value = '  +7.62%  '
$ws.Range('D17').Value = "'19.46"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.89%  '
$ws.Range('D18').Value = '70.571.54'
$ws.Range('E18').Value = '  +7.51%  '
$ws.Range('D19').Value = "'12.63"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.57%  '
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('E21').Value = '  +5.59%  '
$ws.Range('D22').Value = "'492.10"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.45%  '
$ws.Range('D23').Value = "'5.53"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +11.85%  '
$ws.Range('D24').Value = "'16.81"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +16.85%  '
$ws.Range('E25').Value = '  +8.90%  '
$ws.Range('D26').Value = "'90.73"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.06%  '
$ws.Range('E27').Value = '  +6.20%  '
$ws.Range('D28').Value = "'11.22"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.77%  '
$ws.Range('D29').Value = "'9.41"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.20%  '
$ws.Range('D30').Value = "'32.41"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.24%  '
$ws.Range('D31').Value = "'7.74"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +17.66%  '
$ws.Range('D32').Value = "'12.25"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.22%  '
$ws.Range('D33').Value = "'611.55"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.25%  '
$ws.Range('D34').Value = "'65.49"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.63%  '
$ws.Range('E35').Value = '  +7.74%  '
$ws.Range('D36').Value = '0.0₃0833'
$ws.Range('E36').Value = '  +12.38%  '
$ws.Range('D37').Value = "'0.149"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.64%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').Value = "'38.02"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.67%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = "'1.00"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('E40').Value = '  +7.05%  '
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('D42').Value = '3.364.47'
$ws.Range('E42').Value = '  +8.55%  '
$ws.Range('D43').Value = "'3.07"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.58%  '
$ws.Range('E44').Value = '  +7.18%  '
$ws.Range('E45').Value = '  +8.45%  '
$ws.Range('D46').Value = "'3.41"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +7.85%  '
$ws.Range('E47').Value = '  +3.05%  '
$ws.Range('D48').Value = "'9.16"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.22%  '
$ws.Range('D49').Value = "'3.39"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.93%  '
$ws.Range('D50').Value = "'2.74"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +10.42%  '
$ws.Range('E51').Value = '  -0.01%  '
